# edit.ps1
# Adds three new worksheets to the workbook, matching the commit:
#   1. "CypherOutput_Message"  - a copy of the "Message" sheet (connection info
#                                 + the original Cypher query used for CypherOutput)
#   2. "StatOutput"            - a small 4-column summary of file/sample/case/study counts
#   3. "StatOutput_Message"    - connection info + Cypher query block repeated twice,
#                                 the second occurrence's query cell replaced with the
#                                 new Cypher query (the one that produced StatOutput)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) CypherOutput_Message = copy of Message, appended at the end, renamed
# ---------------------------------------------------------------------------
$messageSheet = $wb.Worksheets.Item("Message")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$messageSheet.Copy($null, $lastSheet)
$cypherOutputMessage = $wb.Worksheets.Item($wb.Worksheets.Count)
$cypherOutputMessage.Name = "CypherOutput_Message"

# ---------------------------------------------------------------------------
# 2) StatOutput - header row + one data row of counts
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOutput = $wb.Worksheets.Add($null, $lastSheet)
$statOutput.Name = "StatOutput"

$statOutput.Range("A1").Value = "number_of_files"
$statOutput.Range("B1").Value = "number_of_sample"
$statOutput.Range("C1").Value = "number_of_cases"
$statOutput.Range("D1").Value = "number_of_study"

# Write the counts as genuine text cells (not numbers) by routing them through
# a TEXT() formula and then pasting the computed values back over themselves -
# this avoids Excel's "looks like a number" auto-conversion of a literal
# string assignment while still leaving the cell format untouched (no new
# style entries get created).
$statOutput.Range("A2").Formula = "=TEXT(1,""0"")"
$statOutput.Range("B2").Formula = "=TEXT(2,""0"")"
$statOutput.Range("C2").Formula = "=TEXT(1,""0"")"
$statOutput.Range("D2").Formula = "=TEXT(1,""0"")"
$statOutput.Calculate()
$statRow = $statOutput.Range("A2:D2")
$statRow.Copy($null)
$statRow.PasteSpecial(-4163, $null, $false, $false)

# ---------------------------------------------------------------------------
# 3) StatOutput_Message - Message block, repeated, second block's Cypher row
#    swapped for the new Cypher query used to build StatOutput
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$messageSheet.Copy($null, $lastSheet)
$statOutputMessage = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOutputMessage.Name = "StatOutput_Message"

$newCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Akita']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

for ($i = 1; $i -le 10; $i++) {
    $srcCell = $messageSheet.Cells.Item($i, 1)
    $statOutputMessage.Cells.Item(10 + $i, 1).Value = $srcCell.Value2
}
$statOutputMessage.Range("A18").Value = $newCypher

# ---------------------------------------------------------------------------
# Restore the original active sheet/selection (sheet 1, "CypherOutput") so the
# workbook-level view state matches the pre-edit file instead of drifting to
# whichever sheet we touched last.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
